$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.143.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.139.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.40%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.138.61"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.476"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.658.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.180.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.38%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.148.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.42%  "
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0735"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "450.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0400"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.878.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.271"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.114"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
